# Auto-generated Excel COM-interop script to apply numeric value updates
# described by the OOXML diff (Cactuar_Profits.xlsx / Sheets workbook).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates_ALC = @{
    "H11" = 114.4
    "I11" = 114.4
    "K11" = 114.4
    "M11" = 25.59999999999999
    "H33" = 191.6923
    "I33" = 232.77777
    "K33" = 232.77777
    "M33" = -3.777770000000004
    "H40" = 31777.555
    "J40" = 24428.428
    "L40" = 24428.428
    "N40" = -24778.428
    "H53" = 624
    "I53" = 589.75
    "J53" = 646.8333
    "K53" = 589.75
    "L53" = 646.8333
    "M53" = 47.25
    "N53" = -1920.8333
    "H74" = 5951
    "I74" = 7100
    "J74" = 5185
    "K74" = 7100
    "L74" = 5185
    "M74" = -6164
    "N74" = -7057
    "H76" = 3699.5
    "I76" = 3019.6
    "K76" = 3019.6
    "M76" = -2704.6
    "H77" = 5951
    "I77" = 7100
    "J77" = 5185
    "K77" = 35500
    "L77" = 25925
    "M77" = -30820
    "N77" = -35285
    "H79" = 3699.5
    "I79" = 3019.6
    "K79" = 3019.6
    "M79" = -1927.6
    "H82" = 7293.4
    "J82" = 10000
    "L82" = 30000
    "N82" = -30812
    "H85" = 7293.4
    "J85" = 10000
    "L85" = 30000
    "N85" = -32808
    "H98" = 1596.7727
    "I98" = 1228.5625
    "J98" = 2578.6667
    "K98" = 1228.5625
    "L98" = 2578.6667
    "M98" = 269.4375
    "N98" = -5574.6667
    "H106" = 41670656
    "I106" = 66669150
    "K106" = 66669150
    "M106" = -66668519
    "H111" = 1318.8182
    "I111" = 1428.3
    "J111" = 224
    "K111" = 4284.9
    "L111" = 672
    "M111" = -1217.9
    "N111" = -6806
    "H112" = 4232
    "J112" = 4366.6665
    "L112" = 13099.9995
    "N112" = -15315.9995
    "H122" = 1596.7727
    "I122" = 1228.5625
    "J122" = 2578.6667
    "K122" = 3685.6875
    "L122" = 7736.000100000001
    "M122" = -1235.6875
    "N122" = -12636.0001
    "H132" = 195747.45
    "I132" = 227841.83
    "K132" = 683525.49
    "M132" = -680995.49
    "H135" = 41030.37
    "I135" = 45845.695
    "J135" = 13342.25
    "K135" = 412611.255
    "L135" = 120080.25
    "M135" = -410076.255
    "N135" = -125150.25
    "H137" = 3846.6924
    "I137" = 2855.2222
    "J137" = 4371.5884
    "K137" = 8565.6666
    "L137" = 13114.7652
    "M137" = -6015.6666
    "N137" = -18214.7652
    "H138" = 5460.7334
    "I138" = 1742.8889
    "J138" = 5967.712
    "K138" = 5228.6667
    "L138" = 17903.136
    "M138" = -88.66669999999976
    "N138" = -28183.136
    "H141" = 5775.8335
    "I141" = 5821.4707
    "K141" = 17464.4121
    "M141" = -12284.4121
}
foreach ($cellRef in $updates_ALC.Keys) {
    $ws.Range($cellRef).Value = $updates_ALC[$cellRef]
}

$ws = $wb.Worksheets.Item("ARM")
$updates_ARM = @{
    "H32" = 2985.6418
    "I32" = 1583.2
    "J32" = 9413.5
    "K32" = 1583.2
    "L32" = 9413.5
    "M32" = -1296.2
    "N32" = -9987.5
    "H45" = 3073.875
    "I45" = 5499
    "J45" = 2727.4285
    "K45" = 5499
    "L45" = 2727.4285
    "M45" = -5122
    "N45" = -3481.4285
    "H61" = 14065.324
    "I61" = 7924.9165
    "K61" = 7924.9165
    "M61" = -7712.9165
    "H63" = 129136.75
    "I63" = 170582.33
    "K63" = 170582.33
    "M63" = -169896.33
    "H66" = 129136.75
    "I66" = 170582.33
    "K66" = 852911.6499999999
    "M66" = -849479.6499999999
    "H74" = 3357.1304
    "I74" = 1364.2667
    "K74" = 1364.2667
    "M74" = -490.2666999999999
    "H77" = 3357.1304
    "I77" = 1364.2667
    "K77" = 6821.3335
    "M77" = -2453.3335
    "H102" = 2283843.5
    "I102" = 2283843.5
    "J102" = 0
    "K102" = 2283843.5
    "L102" = 0
    "M102" = -2282221.5
    "H110" = 1025019.1
    "I110" = 1702508.5
    "J110" = 8785
    "K110" = 1702508.5
    "L110" = 8785
    "M110" = -1700463.5
    "N110" = -12875
    "H122" = 349594.06
    "I122" = 692088.9
    "J122" = 7099.25
    "K122" = 2076266.7
    "L122" = 21297.75
    "M122" = -2073816.7
    "N122" = -26197.75
    "H131" = 94998.89999999999
    "J131" = 94998.89999999999
    "L131" = 94998.89999999999
    "N131" = -105078.9
    "H132" = 28091.709
    "I132" = 29637.95
    "J132" = 20360.5
    "K132" = 88913.85000000001
    "L132" = 61081.5
    "M132" = -86383.85000000001
    "N132" = -66141.5
    "H135" = 49999.668
    "J135" = 49999.668
    "L135" = 49999.668
    "N135" = -60139.668
    "H136" = 14065.324
    "I136" = 7924.9165
    "K136" = 23774.7495
    "M136" = -21224.7495
}
foreach ($cellRef in $updates_ARM.Keys) {
    $ws.Range($cellRef).Value = $updates_ARM[$cellRef]
}
$deletions_ARM = @("N102")
foreach ($cellRef in $deletions_ARM) {
    $ws.Range($cellRef).ClearContents()
}

$ws = $wb.Worksheets.Item("BSM")
$updates_BSM = @{
    "H94" = 653106.9399999999
    "I94" = 1142452.9
    "K94" = 1142452.9
    "M94" = -1142001.9
    "H99" = 2085087.9
    "I99" = 2085087.9
    "K99" = 2085087.9
    "M99" = -2083589.9
    "H105" = 200008000
    "J105" = 5999.5
    "L105" = 5999.5
    "N105" = -9493.5
    "H134" = 3672.5
    "I134" = 2284.9656
    "J134" = 6355.067
    "K134" = 6854.8968
    "L134" = 19065.201
    "M134" = -4319.8968
    "N134" = -24135.201
    "H141" = 73217.75
    "J141" = 73217.75
    "L141" = 73217.75
    "N141" = -83577.75
}
foreach ($cellRef in $updates_BSM.Keys) {
    $ws.Range($cellRef).Value = $updates_BSM[$cellRef]
}

$ws = $wb.Worksheets.Item("CRP")
$updates_CRP = @{
    "H16" = 1568.5
    "I16" = 1662.2
    "K16" = 1662.2
    "M16" = -1375.2
    "H31" = 21281006
    "I31" = 40002310
    "J31" = 6795.3184
    "K31" = 40002310
    "L31" = 6795.3184
    "M31" = -40002015
    "N31" = -7385.3184
    "H34" = 21281006
    "I34" = 40002310
    "J34" = 6795.3184
    "K34" = 40002310
    "L34" = 6795.3184
    "M34" = -40002108
    "N34" = -7199.3184
    "H58" = 1003601.1
    "I58" = 1252899.6
    "K58" = 1252899.6
    "M58" = -1252696.6
    "H86" = 4389.5
    "I86" = 3500
    "J86" = 4834.25
    "K86" = 3500
    "L86" = 4834.25
    "M86" = -2377
    "N86" = -7080.25
    "H89" = 4389.5
    "I89" = 3500
    "J89" = 4834.25
    "K89" = 17500
    "L89" = 24171.25
    "M89" = -11884
    "N89" = -35403.25
    "H103" = 34999.5
    "I103" = 13704.4
    "K103" = 13704.4
    "M103" = -12532.4
    "H113" = 1568.5
    "I113" = 1662.2
    "K113" = 1662.2
    "M113" = 507.8
    "H132" = 18523030
    "I132" = 21279392
    "K132" = 63838176
    "M132" = -63835646
    "H136" = 1003601.1
    "I136" = 1252899.6
    "K136" = 3758698.8
    "M136" = -3756148.8
}
foreach ($cellRef in $updates_CRP.Keys) {
    $ws.Range($cellRef).Value = $updates_CRP[$cellRef]
}

$ws = $wb.Worksheets.Item("CUL")
$updates_CUL = @{
    "H3" = 7517.375
    "I3" = 1689.8334
    "K3" = 5069.5002
    "M3" = -4957.5002
    "H14" = 300.16666
    "I14" = 300.16666
    "K14" = 900.4999799999999
    "M14" = -727.4999799999999
    "H26" = 788.8
    "I26" = 236
    "J26" = 3000
    "K26" = 708
    "L26" = 9000
    "M26" = -420
    "N26" = -9576
    "H75" = 4431.1665
    "J75" = 4367.4
    "L75" = 13102.2
    "N75" = -15098.2
    "H78" = 4431.1665
    "J78" = 4367.4
    "L78" = 39306.6
    "N78" = -49290.6
    "H92" = 1002.44446
    "I92" = 879.4
    "K92" = 2638.2
    "M92" = -1390.2
    "H98" = 1506.3
    "I98" = 824.5
    "J98" = 1676.75
    "K98" = 2473.5
    "L98" = 5030.25
    "M98" = -975.5
    "N98" = -8026.25
    "H99" = 3599
    "I99" = 1665
    "J99" = 6500
    "K99" = 4995
    "L99" = 19500
    "M99" = -2749
    "N99" = -23992
    "H116" = 2859.2
    "I116" = 1682.1666
    "J116" = 4624.75
    "K116" = 5046.4998
    "L116" = 13874.25
    "M116" = -1604.4998
    "N116" = -20758.25
    "H128" = 344258.62
    "I128" = 344258.62
    "K128" = 1032775.86
    "M128" = -1027795.86
    "H131" = 14744027
    "J131" = 6631023
    "L131" = 19893069
    "N131" = -19903149
    "H133" = 6834.222
    "J133" = 9196.799999999999
    "L133" = 27590.4
    "N133" = -37710.39999999999
    "H137" = 33083848
    "I137" = 37501750
    "J137" = 18357500
    "K137" = 112505250
    "L137" = 55072500
    "M137" = -112500150
    "N137" = -55082700
    "H138" = 3294
    "I138" = 3617.5
    "J138" = 2000
    "K138" = 10852.5
    "L138" = 6000
    "M138" = -5712.5
    "N138" = -16280
}
foreach ($cellRef in $updates_CUL.Keys) {
    $ws.Range($cellRef).Value = $updates_CUL[$cellRef]
}

$ws = $wb.Worksheets.Item("GSM")
$updates_GSM = @{
    "H39" = 110500
    "J39" = 110500
    "L39" = 110500
    "N39" = -111564
    "H70" = 1366746.2
    "J70" = 7542.4165
    "L70" = 7542.4165
    "N70" = -8082.4165
    "H73" = 1366746.2
    "J73" = 7542.4165
    "L73" = 7542.4165
    "N73" = -9414.416499999999
    "H80" = 1119152.6
    "I80" = 1517208.2
    "J80" = 24499.5
    "K80" = 1517208.2
    "L80" = 24499.5
    "M80" = -1516210.2
    "N80" = -26495.5
    "H83" = 1119152.6
    "I83" = 1517208.2
    "J83" = 24499.5
    "K83" = 7586041
    "L83" = 122497.5
    "M83" = -7581049
    "N83" = -132481.5
    "H97" = 1679.2285
    "I97" = 1355.6923
    "J97" = 2613.889
    "K97" = 1355.6923
    "L97" = 2613.889
    "M97" = -859.6922999999999
    "N97" = -3605.889
    "H122" = 293611.4
    "I122" = 582585.9
    "J122" = 4636.9473
    "K122" = 1747757.7
    "L122" = 13910.8419
    "M122" = -1745307.7
    "N122" = -18810.8419
    "H126" = 3395.7585
    "I126" = 2012.0869
    "K126" = 6036.2607
    "M126" = -3566.2607
    "H132" = 3686.6226
    "I132" = 3838.814
    "J132" = 3032.2
    "K132" = 11516.442
    "L132" = 9096.599999999999
    "M132" = -8986.441999999999
    "N132" = -14156.6
    "H135" = 105000
    "J135" = 120000
    "L135" = 120000
    "N135" = -130140
    "H136" = 0
    "J136" = 0
    "L136" = 0
}
foreach ($cellRef in $updates_GSM.Keys) {
    $ws.Range($cellRef).Value = $updates_GSM[$cellRef]
}
$deletions_GSM = @("N136")
foreach ($cellRef in $deletions_GSM) {
    $ws.Range($cellRef).ClearContents()
}

$ws = $wb.Worksheets.Item("LTW")
$updates_LTW = @{
    "H7" = 4114.2905
    "J7" = 5399.8667
    "L7" = 5399.8667
    "N7" = -5623.8667
    "H22" = 1153.2858
    "I22" = 1632.8334
    "K22" = 1632.8334
    "M22" = -1337.8334
    "H27" = 1153.2858
    "I27" = 1632.8334
    "K27" = 1632.8334
    "M27" = -1525.8334
    "H55" = 282.34784
    "I55" = 445.83334
    "J55" = 104
    "K55" = 445.83334
    "L55" = 104
    "M55" = -272.83334
    "N55" = -450
    "H68" = 22727272
    "I68" = 22727272
    "K68" = 22727272
    "M68" = -22726523
    "H71" = 22727272
    "I71" = 22727272
    "K71" = 113636360
    "M71" = -113632616
    "H122" = 6464.2144
    "I122" = 5061.75
    "J122" = 8334.166999999999
    "K122" = 15185.25
    "L122" = 25002.501
    "M122" = -12735.25
    "N122" = -29902.501
    "H126" = 4114.2905
    "J126" = 5399.8667
    "L126" = 16199.6001
    "N126" = -21139.6001
    "H132" = 4305.85
    "J132" = 4495
    "L132" = 13485
    "N132" = -18545
    "H136" = 3987.32
    "I136" = 3986.2197
    "J136" = 3998.4443
    "K136" = 11958.6591
    "L136" = 11995.3329
    "M136" = -9408.659100000001
    "N136" = -17095.3329
}
foreach ($cellRef in $updates_LTW.Keys) {
    $ws.Range($cellRef).Value = $updates_LTW[$cellRef]
}

$ws = $wb.Worksheets.Item("WVR")
$updates_WVR = @{
    "H115" = 18730
    "I115" = 4912.5
    "J115" = 74000
    "K115" = 4912.5
    "L115" = 74000
    "M115" = -3345.5
    "N115" = -77134
    "H126" = 2380.7222
    "I126" = 1621
    "J126" = 5039.75
    "K126" = 4863
    "L126" = 15119.25
    "M126" = -2393
    "N126" = -20059.25
    "H132" = 1858518.1
    "I132" = 2653490.5
    "J132" = 3582.6667
    "K132" = 7960471.5
    "L132" = 10748.0001
    "M132" = -7957941.5
    "N132" = -15808.0001
    "H136" = 6933.05
    "I136" = 4641.6
    "J136" = 9224.5
    "K136" = 13924.8
    "L136" = 27673.5
    "M136" = -11374.8
    "N136" = -32773.5
}
foreach ($cellRef in $updates_WVR.Keys) {
    $ws.Range($cellRef).Value = $updates_WVR[$cellRef]
}

Write-Output "Applied all Sheets updates successfully."